$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1200.42
$wsSummary.Range("B4").Value = 0.42
$wsSummary.Range("B5").Value = 0.16
$wsSummary.Range("B6").Value = 53
$wsSummary.Range("B8").Value = 29
$wsSummary.Range("B9").Value = 28.3

# --- Strategy Status sheet ---
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C4").Value = 100.42
$wsStatus.Range("D4").Value = 53
$wsStatus.Range("E4").Value = 0.42
$wsStatus.Range("F4").Value = 0.42
$wsStatus.Range("G4").Value = 28.3

# --- New trade row (#53) to append to both "All Trades" and "MarketMaking" sheets ---
$newRow = @{
    A = 53
    B = "2026-02-17"
    C = "15:42:54"
    D = "MarketMaking"
    E = "UP"
    F = 0.33
    G = 0.23
    H = "CLOSED"
    I = -30.303
    J = -0.1
    K = 100.42
    L = 0
    M = 0
    N = 0.6
    O = "Normal spread capture: 19600 bps"
    P = "early_exit"
    Q = 0.15
}

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("A54").Value = $newRow.A
    # "2026-02-17" looks like a date, so Excel would normally auto-convert it
    # to a date serial number. Force the cell to text first so it is stored
    # the same way as the existing date cells in this column (inline/shared
    # string), then restore the default "Normal" style so no stray number
    # format lingers on the cell.
    $ws.Range("B54").NumberFormat = "@"
    $ws.Range("B54").Value = $newRow.B
    $ws.Range("B54").Style = "Normal"
    $ws.Range("C54").Value = $newRow.C
    $ws.Range("D54").Value = $newRow.D
    $ws.Range("E54").Value = $newRow.E
    $ws.Range("F54").Value = $newRow.F
    $ws.Range("G54").Value = $newRow.G
    $ws.Range("H54").Value = $newRow.H
    $ws.Range("I54").Value = $newRow.I
    $ws.Range("J54").Value = $newRow.J
    $ws.Range("K54").Value = $newRow.K
    $ws.Range("L54").Value = $newRow.L
    $ws.Range("M54").Value = $newRow.M
    $ws.Range("N54").Value = $newRow.N
    $ws.Range("O54").Value = $newRow.O
    $ws.Range("P54").Value = $newRow.P
    $ws.Range("Q54").Value = $newRow.Q
}
